$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was a duplicate of row 1 (Jessie Marlowe). Replace it with a new,
# distinct person: Michelle Norton.
$ws.Range("A2").Value = "Michelle"
$ws.Range("B2").Value = "Norton"
$ws.Range("C2").Value = "Aperture Inc."
$ws.Range("D2").Value = "Scientist"
$ws.Range("E2").Value = "13 White Rabbit Street"
$ws.Range("F2").Value = "mnorton@aperture.us"
$ws.Range("G2").Value = 40731254562

# Row 5 was a duplicate of row 4 (Michael Robertson). Replace it with a new,
# distinct person: Jane Dorsey.
$ws.Range("A5").Value = "Jane"
$ws.Range("B5").Value = "Dorsey"
$ws.Range("C5").Value = "MediCare"
$ws.Range("D5").Value = "Medical Engineer"
$ws.Range("E5").Value = "11 Crown Street"
$ws.Range("F5").Value = "jdorsey@mc.com"
$ws.Range("G5").Value = 40791345621

# Row 9 was a duplicate of row 8 (Doug Derrick). Replace it with a new,
# distinct person: Lara Palmer.
$ws.Range("A9").Value = "Lara"
$ws.Range("B9").Value = "Palmer"
$ws.Range("C9").Value = "Timepath Inc."
$ws.Range("D9").Value = "Programmer"
$ws.Range("E9").Value = "87 Orange Street"
$ws.Range("F9").Value = "lpalmer@timepath.co.uk"
$ws.Range("G9").Value = 40731653845
